# Apply the "Add UML sequence diagrams content" edit to slide 1:
#  - remove the (empty, unused) Title placeholder from the diagram slide
#  - slide the remaining UML notation diagram shapes up to reclaim the
#    vertical space that the title placeholder used to occupy

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Drop the custom "first slide number" override (back to the default
# start-at-1 numbering).
$p.PageSetup.FirstSlideNumber = 1

# ppPlaceholderTitle = 1 ; msoPlaceholder = 14
$ppPlaceholderTitle = 1
$msoPlaceholder = 14

# 1) Delete the empty title placeholder shape.
for ($i = $s.Shapes.Count; $i -ge 1; $i--) {
    $sh = $s.Shapes.Item($i)
    if ($sh.Type -eq $msoPlaceholder) {
        if ($sh.PlaceholderFormat.Type -eq $ppPlaceholderTitle) {
            $sh.Delete()
        }
    }
}

# 2) Move every remaining shape up by 1895386 EMU (≈149.243 pt) to close
#    the gap left behind by the removed title. Target tops are given in
#    EMU (native OOXML units, 914400 EMU/in, 12700 EMU/pt) and converted
#    to points for the .Top property; a half-EMU nudge keeps the
#    point<->EMU float round-trip from truncating one unit short.
$emuPerPt = 12700.0
$halfEmuPt = 0.5 / $emuPerPt

function Set-TopEmu($shape, [double]$targetEmu) {
    $shape.Top = ($targetEmu / $emuPerPt) + $halfEmuPt
}

$targets = @{
    "Rectangle 16"       = 2993082.0
    "Rectangle 17"       = 3356580.0
    "Folded Corner 18"   = 2448014.0
    "TextBox 19"         = 2448015.0
    "Rectangle 20"       = 2981414.0
    "Rectangle 21"       = 3356580.0
    "Elbow Connector 22" = 3541246.0
    "Rectangle 28"       = 3667214.0
    "Rectangle 29"       = 3667214.0
}

for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $sh = $s.Shapes.Item($i)
    if ($targets.ContainsKey($sh.Name)) {
        Set-TopEmu $sh $targets[$sh.Name]
    }
}
